# Bump the two "INC" suffix counters used to build the test usernames/emails.
# Sheet "Sheet1" holds the formulas (CONCATENATE($G,$I$2|$I$23)) that generate
# the usernames; bumping I2 (10->11) and I23 (15->16) recalculates columns
# A:C there. The "login" and "order" sheets hold static copies of those same
# generated strings (no formulas), so they are updated in lock-step here.

$wb   = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Sheet1")
$login = $wb.Worksheets.Item("login")
$order = $wb.Worksheets.Item("order")

$oldSuffix1 = "10"
$newSuffix1 = "11"
$oldSuffix2 = "15"
$newSuffix2 = "16"

# --- Sheet1: bump the two source counters; formulas recalc automatically ---
$data.Range("I2").Value = [int]$newSuffix1
$data.Range("I23").Value = [int]$newSuffix2

# --- login sheet: rows 2-21 mirror the first block (EthanBaker10, ...) ---
for ($r = 2; $r -le 21; $r++) {
    $oldG = $login.Range("G$r").Value()
    $newG = $oldG -replace "$oldSuffix1$", $newSuffix1
    $login.Range("G$r").Value = $newG
    $login.Range("H$r").Value = $newG

    $oldI = $login.Range("I$r").Value()
    $newI = $oldI -replace "$oldSuffix1@", "$newSuffix1@"
    $login.Range("I$r").Value = $newI
}

# --- order sheet: rows 2-21 mirror the second block (DonnellJernigan15, ...) ---
for ($r = 2; $r -le 21; $r++) {
    $oldR = $order.Range("R$r").Value()
    $newR = $oldR -replace "$oldSuffix2$", $newSuffix2
    $order.Range("R$r").Value = $newR
    $order.Range("S$r").Value = $newR

    $oldT = $order.Range("T$r").Value()
    $newT = $oldT -replace "$oldSuffix2@", "$newSuffix2@"
    $order.Range("T$r").Value = $newT
}

# --- Cosmetic view state (best-effort; window geometry/topLeftCell are not
#     part of the workbook's data model so may not round-trip, but set them
#     in case the host honors it) ---
$data.Activate()
$excel.Goto($data.Range("A19"), $false)
$data.Range("A23:C42").Select()

$win = $excel.ActiveWindow
$win.Top = 1500
$win.Left = 0
$win.Width = 23040
$win.Height = 12180
